# Fruta / hortaliza, semanal
# Inserts a new weekly record at row 21 (shifting the existing rows 21-78
# down to 22-79) on the "Poroto granado" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 21..78 down by one to make room for the new record.
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new weekly observation.
$ws.Range("A21").Value = 4
$ws.Range("B21").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C21").Value = "Los Lagos"
$ws.Range("D21").Value = 45044
$ws.Range("D21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = 100112030
$ws.Range("G21").Value = "Poroto granado"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 50
$ws.Range("K21").Value = 36000
$ws.Range("L21").Value = 36000
$ws.Range("M21").Value = 36000
$ws.Range("N21").Value = "$/saco 25 kilos"
$ws.Range("O21").Value = "Región Metropolitana"
$ws.Range("P21").Value = 1440
$ws.Range("Q21").Value = 25
$ws.Range("R21").Value = "Hortaliza"
